$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: new header labels in Q:U (order matters for shared string table) ---
$ws.Range("R30").Value = "Mult for minute"
$ws.Range("S30").Value = "Games per minute"
$ws.Range("Q30").Value = "Games played"
$ws.Range("T30").Value = "Desired Minutes"
$ws.Range("U30").Value = "Required Games"
$ws.Range("Q30:U30").WrapText = $true
$ws.Rows.Item(30).RowHeight = 39

# --- Remove old V:W "runs per depth" helper columns (rows 31-41) ---
$ws.Range("V31:W41").ClearContents()

# --- New T (desired minutes) / U (required games) columns, rows 31-39 ---
$ws.Range("T31").Value = 1
$ws.Range("U31").Formula = "=S31*T31"

$ws.Range("T32").Value = 3
$ws.Range("U32").Formula = "=S32*T32"

$ws.Range("T33").Value = 5
$ws.Range("U33").Formula = "=S33*T33"

$ws.Range("T34").Value = 25
$ws.Range("U34").Formula = "=S34*T34"

$ws.Range("T35").Value = 60
$ws.Range("T36").Value = 60
$ws.Range("T37").Value = 120
$ws.Range("T38").Value = 360
$ws.Range("T39").Value = 720
$ws.Range("U35:U39").Formula = "=S35*T35"

$ws.Range("U31:U39").NumberFormat = "#,##0"

# --- New summary cell ---
$ws.Range("T44").Formula = "=SUM(T31:T39)/60"

# --- Selection state ---
$ws.Range("T40").Select() | Out-Null
